$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: shift E5/G5 content right, add new H5 text ---
# Old: E5 = "Must be integrated within a web application. "
#      G5 = "Additional resource options if the student does not understand the steps."
#      H5 = (empty)
# New: E5 = "Must be integrated within a web application. "   (unchanged content, new shared-string slot)
#      F5 = "Must Have"                                         (unchanged)
#      G5 = "Additional resource options if the student does not understand the steps."
#      H5 = "Happy brithday playing whenever you git the calculate button."
$ws.Range("E5").Value = "Must be integrated within a web application. "
$ws.Range("G5").Value = "Additional resource options if the student does not understand the steps."
$ws.Range("H5").Value = "Happy brithday playing whenever you git the calculate button."

# --- Row 9: replace the user story text (A9) and add a new Should Have (H9) ---
$ws.Range("A9").Value = "I Sam, as a user who must take Physics to obtain my degree have trouble because the fundementals of Physics can be  frusterating for students who tend to be visual learners. This web application is meant to alleviate those troubles by adding graphics to the calculations to help students visualize the concepts by seeing graphics of how interactions of objects occur and the different forces/energies that they interact with as well as the objects themselves. "
$ws.Range("H9").Value = "An email or credit card to access the application."

# --- Row 10: populate the previously-empty Must/Should/Could/Won't Have cells, and grow the row height ---
$ws.Range("E10").Value = "Factually correct information pertaining to physics fundementals displayed during runtime on the application."
$ws.Range("F10").Value = "Easy to understand UI interface."
$ws.Range("G10").Value = "A way to save specific results to come back to the same version of a problem at a later time without having to reset. Cookie based system stored locally. "
$ws.Range("H10").Value = "Use Cookies to do anything malicious on end users computer. "
$ws.Rows.Item(10).RowHeight = 128

# --- Update the visible selection to match the saved view state ---
$ws.Range("F10").Select()
